# Auto-generated Excel COM-interop script to apply the Zalera_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(103, 8).Value = 1121.6154
$ws.Cells.Item(103, 9).Value = 1395.2858
$ws.Cells.Item(103, 10).Value = 802.3333
$ws.Cells.Item(103, 11).Value = 4185.857400000001
$ws.Cells.Item(103, 12).Value = 2406.9999
$ws.Cells.Item(103, 13).Value = -3599.857400000001
$ws.Cells.Item(103, 14).Value = -3578.9999
$ws.Cells.Item(113, 8).Value = 1000006
$ws.Cells.Item(113, 10).Value = 1000006
$ws.Cells.Item(113, 12).Value = 1000006
$ws.Cells.Item(113, 14).Value = -1006514
$ws.Cells.Item(136, 8).Value = 99773
$ws.Cells.Item(136, 10).Value = 99773
$ws.Cells.Item(136, 12).Value = 99773
$ws.Cells.Item(136, 14).Value = -109973
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(6, 8).Value = 1006001
$ws.Cells.Item(6, 9).Value = 1006001
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 1006001
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = -1005828
$ws.Cells.Item(6, 14).Value = $null
$ws.Cells.Item(53, 8).Value = 5000
$ws.Cells.Item(53, 9).Value = 5000
$ws.Cells.Item(53, 11).Value = 5000
$ws.Cells.Item(53, 13).Value = -4318
$ws.Cells.Item(110, 8).Value = 41667864
$ws.Cells.Item(110, 9).Value = 50000740
$ws.Cells.Item(110, 11).Value = 50000740
$ws.Cells.Item(110, 13).Value = -49998695
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(36, 8).Value = 16018.5
$ws.Cells.Item(36, 9).Value = 16018.5
$ws.Cells.Item(36, 11).Value = 16018.5
$ws.Cells.Item(36, 13).Value = -15484.5
$ws.Cells.Item(46, 8).Value = 20339.285
$ws.Cells.Item(46, 10).Value = 20339.285
$ws.Cells.Item(46, 12).Value = 20339.285
$ws.Cells.Item(46, 14).Value = -20935.285
$ws.Cells.Item(57, 8).Value = 57950.332
$ws.Cells.Item(57, 10).Value = 66398.60000000001
$ws.Cells.Item(57, 12).Value = 66398.60000000001
$ws.Cells.Item(57, 14).Value = -67838.60000000001
$ws.Cells.Item(80, 8).Value = 74928.28999999999
$ws.Cells.Item(80, 9).Value = 704.6667
$ws.Cells.Item(80, 10).Value = 130596
$ws.Cells.Item(80, 11).Value = 704.6667
$ws.Cells.Item(80, 12).Value = 130596
$ws.Cells.Item(80, 13).Value = 293.3333
$ws.Cells.Item(80, 14).Value = -132592
$ws.Cells.Item(83, 8).Value = 74928.28999999999
$ws.Cells.Item(83, 9).Value = 704.6667
$ws.Cells.Item(83, 10).Value = 130596
$ws.Cells.Item(83, 11).Value = 3523.3335
$ws.Cells.Item(83, 12).Value = 652980
$ws.Cells.Item(83, 13).Value = 1468.6665
$ws.Cells.Item(83, 14).Value = -662964
$ws.Cells.Item(86, 8).Value = 501699.66
$ws.Cells.Item(86, 9).Value = 2549.5
$ws.Cells.Item(86, 11).Value = 2549.5
$ws.Cells.Item(86, 13).Value = -1426.5
$ws.Cells.Item(89, 8).Value = 501699.66
$ws.Cells.Item(89, 9).Value = 2549.5
$ws.Cells.Item(89, 11).Value = 12747.5
$ws.Cells.Item(89, 13).Value = -7131.5
$ws.Cells.Item(107, 8).Value = 1439.2413
$ws.Cells.Item(107, 9).Value = 1493.7037
$ws.Cells.Item(107, 10).Value = 704
$ws.Cells.Item(107, 11).Value = 1493.7037
$ws.Cells.Item(107, 12).Value = 704
$ws.Cells.Item(107, 13).Value = 426.2963
$ws.Cells.Item(107, 14).Value = -4544
$ws.Cells.Item(134, 8).Value = 19282.363
$ws.Cells.Item(134, 9).Value = 33848.25
$ws.Cells.Item(134, 11).Value = 101544.75
$ws.Cells.Item(134, 13).Value = -99009.75
$ws.Cells.Item(136, 8).Value = 57950.332
$ws.Cells.Item(136, 10).Value = 66398.60000000001
$ws.Cells.Item(136, 12).Value = 66398.60000000001
$ws.Cells.Item(136, 14).Value = -76598.60000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 4815.5
$ws.Cells.Item(16, 9).Value = 3801.8333
$ws.Cells.Item(16, 11).Value = 3801.8333
$ws.Cells.Item(16, 13).Value = -3514.8333
$ws.Cells.Item(33, 8).Value = 1510.3334
$ws.Cells.Item(33, 9).Value = 1510.3334
$ws.Cells.Item(33, 11).Value = 1510.3334
$ws.Cells.Item(33, 13).Value = -1131.3334
$ws.Cells.Item(39, 8).Value = 5860
$ws.Cells.Item(39, 9).Value = 5860
$ws.Cells.Item(39, 11).Value = 5860
$ws.Cells.Item(39, 13).Value = -5469
$ws.Cells.Item(47, 8).Value = 10000
$ws.Cells.Item(47, 9).Value = 10000
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 11).Value = 10000
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 13).Value = -9434
$ws.Cells.Item(47, 14).Value = $null
$ws.Cells.Item(49, 8).Value = 5860
$ws.Cells.Item(49, 9).Value = 5860
$ws.Cells.Item(49, 11).Value = 5860
$ws.Cells.Item(49, 13).Value = -5678
$ws.Cells.Item(107, 8).Value = 628.6799999999999
$ws.Cells.Item(107, 9).Value = 652.913
$ws.Cells.Item(107, 11).Value = 652.913
$ws.Cells.Item(107, 13).Value = 1267.087
$ws.Cells.Item(113, 8).Value = 4815.5
$ws.Cells.Item(113, 9).Value = 3801.8333
$ws.Cells.Item(113, 11).Value = 3801.8333
$ws.Cells.Item(113, 13).Value = -1631.8333
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 1389.8
$ws.Cells.Item(113, 10).Value = 1316.3334
$ws.Cells.Item(113, 12).Value = 3949.0002
$ws.Cells.Item(113, 14).Value = -8289.0002
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 12806.322
$ws.Cells.Item(70, 9).Value = 11598.218
$ws.Cells.Item(70, 11).Value = 11598.218
$ws.Cells.Item(70, 13).Value = -11328.218
$ws.Cells.Item(73, 8).Value = 12806.322
$ws.Cells.Item(73, 9).Value = 11598.218
$ws.Cells.Item(73, 11).Value = 11598.218
$ws.Cells.Item(73, 13).Value = -10662.218
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 14).Value = $null
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 14).Value = $null
$ws.Cells.Item(109, 8).Value = 20285
$ws.Cells.Item(109, 10).Value = 20285
$ws.Cells.Item(109, 12).Value = 20285
$ws.Cells.Item(109, 14).Value = -22365
$ws.Cells.Item(126, 8).Value = 3721
$ws.Cells.Item(126, 9).Value = 2892.5715
$ws.Cells.Item(126, 10).Value = 4445.875
$ws.Cells.Item(126, 11).Value = 8677.7145
$ws.Cells.Item(126, 12).Value = 13337.625
$ws.Cells.Item(126, 13).Value = -6207.7145
$ws.Cells.Item(126, 14).Value = -18277.625
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 335933
$ws.Cells.Item(7, 9).Value = 501399.5
$ws.Cells.Item(7, 10).Value = 5000
$ws.Cells.Item(7, 11).Value = 501399.5
$ws.Cells.Item(7, 12).Value = 5000
$ws.Cells.Item(7, 13).Value = -501287.5
$ws.Cells.Item(7, 14).Value = -5224
$ws.Cells.Item(53, 8).Value = 15000
$ws.Cells.Item(53, 10).Value = 15000
$ws.Cells.Item(53, 12).Value = 15000
$ws.Cells.Item(53, 14).Value = -16036
$ws.Cells.Item(55, 8).Value = 506.69232
$ws.Cells.Item(55, 9).Value = 570.3
$ws.Cells.Item(55, 10).Value = 294.66666
$ws.Cells.Item(55, 11).Value = 570.3
$ws.Cells.Item(55, 12).Value = 294.66666
$ws.Cells.Item(55, 13).Value = -397.3
$ws.Cells.Item(55, 14).Value = -640.66666
$ws.Cells.Item(61, 8).Value = 9334.333000000001
$ws.Cells.Item(61, 9).Value = 9334.333000000001
$ws.Cells.Item(61, 11).Value = 9334.333000000001
$ws.Cells.Item(61, 13).Value = -9132.333000000001
$ws.Cells.Item(68, 8).Value = 2386
$ws.Cells.Item(68, 9).Value = 2546.6667
$ws.Cells.Item(68, 11).Value = 2546.6667
$ws.Cells.Item(68, 13).Value = -1797.6667
$ws.Cells.Item(71, 8).Value = 2386
$ws.Cells.Item(71, 9).Value = 2546.6667
$ws.Cells.Item(71, 11).Value = 12733.3335
$ws.Cells.Item(71, 13).Value = -8989.333500000001
$ws.Cells.Item(113, 8).Value = 9334.333000000001
$ws.Cells.Item(113, 9).Value = 9334.333000000001
$ws.Cells.Item(113, 11).Value = 9334.333000000001
$ws.Cells.Item(113, 13).Value = -7164.333000000001
$ws.Cells.Item(116, 8).Value = 77840
$ws.Cells.Item(116, 10).Value = 77840
$ws.Cells.Item(116, 12).Value = 77840
$ws.Cells.Item(116, 14).Value = -87018
$ws.Cells.Item(122, 8).Value = 6826.273
$ws.Cells.Item(122, 9).Value = 5376.6665
$ws.Cells.Item(122, 11).Value = 16129.9995
$ws.Cells.Item(122, 13).Value = -13679.9995
$ws.Cells.Item(126, 8).Value = 335933
$ws.Cells.Item(126, 9).Value = 501399.5
$ws.Cells.Item(126, 10).Value = 5000
$ws.Cells.Item(126, 11).Value = 1504198.5
$ws.Cells.Item(126, 12).Value = 15000
$ws.Cells.Item(126, 13).Value = -1501728.5
$ws.Cells.Item(126, 14).Value = -19940
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(52, 8).Value = 12277.286
$ws.Cells.Item(52, 10).Value = 39899
$ws.Cells.Item(52, 12).Value = 39899
$ws.Cells.Item(52, 14).Value = -40351
$ws.Cells.Item(107, 8).Value = 2654.25
$ws.Cells.Item(107, 9).Value = 2149.5715
$ws.Cells.Item(107, 11).Value = 6448.7145
$ws.Cells.Item(107, 13).Value = -4528.7145
$ws.Cells.Item(135, 8).Value = 74357
$ws.Cells.Item(135, 10).Value = 74357
$ws.Cells.Item(135, 12).Value = 74357
$ws.Cells.Item(135, 14).Value = -84497
$ws.Cells.Item(136, 8).Value = 4086.7188
$ws.Cells.Item(136, 9).Value = 2950.875
$ws.Cells.Item(136, 11).Value = 8852.625
$ws.Cells.Item(136, 13).Value = -6302.625

Write-Host "Applied 206 cell changes across 8 sheets"